# Update cryptocurrency price/volume snapshot cells on Sheet1.
# Column D ("Price") and column E ("Volume(1h)") are stored as literal text
# (not numbers) in the workbook, so values that look numeric (e.g. "1.00",
# "0.999") must be forced to stay text -- otherwise Excel/COM auto-converts
# them to real numbers and strips formatting like trailing/leading zeros.
# We do that by temporarily switching the cell to the Text number format,
# assigning the literal string, then restoring the default "Normal" style so
# the cell's style index is unaffected by the round-trip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.562.12'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '3.115.19'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = '3.113.89'
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.471'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.410'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.39%  '
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("D14").Value = '3.648.10'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("D17").Value = '57.680.35'
$ws.Range("D18").Value = '3.124.34'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '360.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '0.0₃0865'
$ws.Range("E28").Value = '  -3.16%  '
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.17%  '
$ws.Range("E33").Value = '  +4.00%  '
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("E39").Value = '  +3.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0668'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("D41").Value = '2.485.09'
$ws.Range("E41").Value = '  +5.56%  '
$ws.Range("E42").Value = '  -3.29%  '
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '37.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0269'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.976'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("E51").Value = '  +2.46%  '
